# Quarterly financials update: a new quarter (period ending 2018-09-30,
# serial 43373) is inserted as the newest (leftmost) data column on the
# "CAPR" sheet, shifting the previously existing quarters one column to
# the right (old D:K -> new E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D. This shifts all data currently
#    in D:K one column to the right, into E:L, and keeps the column
#    heading/styles intact.
$ws.Columns("D").Insert()

# 2. The freshly inserted column D has no formatting yet. Copy the
#    formatting (number formats / styles) from the column immediately to
#    its right (E, which holds the most recent pre-existing quarter) so
#    that the new quarter column renders the same way (dates as dates,
#    financial figures with the same number format) across the used
#    range of the sheet.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new quarter's data (column D) for every row that has
#    figures in this report. Blank/header rows are intentionally left
#    untouched (they stay blank, matching their neighboring columns).
$newQuarterValues = @{
    7 = 43373
    8 = 200
    9 = "NA"
    10 = "NA"
    12 = 3100
    13 = 0
    14 = "NA"
    15 = 0
    17 = 4400
    18 = -4200
    20 = 100
    21 = -4100
    22 = 0
    23 = -4100
    24 = 0
    25 = 0
    26 = -4100
    27 = -4100
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -100
    33 = -4100
    34 = 0
    35 = -4100
    38 = 43373
    41 = 4400
    42 = 6000
    43 = 200
    44 = 0
    45 = 800
    46 = 11400
    47 = 0
    48 = 600
    49 = 100
    50 = 0
    51 = 0
    52 = 200
    53 = 0
    54 = 12200
    57 = 1900
    58 = "NA"
    59 = "NA"
    60 = 1900
    61 = 3400
    62 = "NA"
    63 = 0
    64 = 0
    65 = 0
    66 = 5300
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -63500
    73 = 0
    74 = 0
    75 = 0
    76 = 7000
    77 = 0
    80 = 43373
    81 = -4100
    83 = 0
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -3000
    91 = -100
    92 = 0
    93 = 0
    94 = -100
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 1100
    101 = 0
    102 = -2100
}

foreach ($row in $newQuarterValues.Keys) {
    $ws.Range("D$row").Value = $newQuarterValues[$row]
}

# 4. Resize the new column to fit its contents, similar to the other
#    best-fit columns on this sheet.
$ws.Columns("D").AutoFit()
